# Refresh the coin Price (D) / Volume(1h) (E) columns of the crypto
# tracker sheet with the latest scraped figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.453.90'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '1.788.94'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").Value = '''306.24'
$ws.Range("E6").Value = '  -0.42%  '
$ws.Range("D7").Value = '''0.4254'
$ws.Range("E7").Value = '  +1.44%  '
$ws.Range("D8").Value = '''0.3622'
$ws.Range("E8").Value = '  +1.67%  '
$ws.Range("D9").Value = '''0.07157'
$ws.Range("E9").Value = '  +1.05%  '
$ws.Range("D10").Value = '''0.8520'
$ws.Range("E10").Value = '  +0.93%  '
$ws.Range("D11").Value = '''20.58'
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("D12").Value = '1.836.05'
$ws.Range("E12").Value = '  +1.75%  '
$ws.Range("D13").Value = '''6.499'
$ws.Range("E13").Value = '  +2.57%  '
$ws.Range("D14").Value = '''5.270'
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = '''0.06882'
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '''79.80'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").Value = '''0.000008815'
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '''15.07'
$ws.Range("E20").Value = '  +0.49%  '
$ws.Range("D21").Value = '26.480.30'
$ws.Range("E21").Value = '  -1.99%  '
$ws.Range("D22").Value = '''5.141'
$ws.Range("E22").Value = '  +1.72%  '
$ws.Range("D23").Value = '''11.04'
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("D24").Value = '2.037.19'
$ws.Range("E24").Value = '  +1.08%  '
$ws.Range("D25").Value = '''152.06'
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").Value = '''1.818'
$ws.Range("E26").Value = '  -5.95%  '
$ws.Range("D27").Value = '''18.17'
$ws.Range("E27").Value = '  +0.45%  '
$ws.Range("D28").Value = '''5.154'
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("D29").Value = '''1.902'
$ws.Range("E29").Value = '  +16.30%  '
$ws.Range("D30").Value = '''114.68'
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("E31").Value = '  -0.72%  '
$ws.Range("D32").Value = '''0.7451'
$ws.Range("E32").Value = '  +3.97%  '
$ws.Range("D33").Value = '''1.144'
$ws.Range("E33").Value = '  +6.69%  '
$ws.Range("D34").Value = '''4.353'
$ws.Range("E34").Value = '  +1.63%  '
$ws.Range("D35").Value = '''2.758'
$ws.Range("E35").Value = '  -3.56%  '
$ws.Range("D36").Value = '''1.002'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '''1.110'
$ws.Range("E37").Value = '  +3.35%  '
$ws.Range("D38").Value = '''0.05160'
$ws.Range("E38").Value = '  +1.23%  '
$ws.Range("D39").Value = '''0.01898'
$ws.Range("E39").Value = '  -0.11%  '
$ws.Range("E40").Value = '  +0.83%  '
$ws.Range("D41").Value = '''0.1621'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").Value = '''2.612'
$ws.Range("E42").Value = '  +1.66%  '
$ws.Range("D43").Value = '''6.398'
$ws.Range("E43").Value = '  +7.20%  '
$ws.Range("D44").Value = '''8.265'
$ws.Range("E44").Value = '  +3.13%  '
$ws.Range("D45").Value = '''105.45'
$ws.Range("E45").Value = '  +0.98%  '
$ws.Range("D46").Value = '''10.27'
$ws.Range("E46").Value = '  +1.37%  '
$ws.Range("D47").Value = '''1.002'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("D48").Value = '''1.641'
$ws.Range("E48").Value = '  +3.17%  '
$ws.Range("D49").Value = '''0.4522'
$ws.Range("E49").Value = '  +0.46%  '
$ws.Range("D50").Value = '''0.06204'
$ws.Range("E50").Value = '  -1.51%  '
$ws.Range("D51").Value = '''1.772'
$ws.Range("E51").Value = '  +5.83%  '
